$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (date + 12 measurement columns B:M) appended after row 9.
# $null entries correspond to cells left blank (as in the source diff).
$rows = @(
    @(45757, $null, $null, $null, $null, $null, $null, $null, $null, 31.2, 33.299999999999997, 27.7, 34),
    @(45758, 36.299999999999997, 41.3, 38.6, 40.1, 38.299999999999997, 32.9, 31.8, 38.299999999999997, 31, 32.4, 26.7, 33),
    @(45396, 36.299999999999997, 42.1, 38.6, 40.200000000000003, 38.299999999999997, 33.1, 32.6, 37.1, 27.8, 28.5, 24.8, 30)
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    for ($c = 1; $c -le 13; $c++) {
        $val = $data[$c - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# Copy the formatting (number format, borders, style index) from the last
# existing data row (9) down across the new rows (10-12), matching how the
# original rows were styled.
$src = $ws.Range("A9:M9")
$dst = $ws.Range("A10:M12")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D17").Select()
